# Adding more docs and training data: new Q&A rows (11-18) appended to the
# Message / Bot Response table on Sheet1, plus formatting cleanup of the
# existing rows 6-10 (drop one-off fonts/row heights back to the sheet's
# normal style) and a couple of minor view/layout tweaks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Copy the "normal" style (style index 2, used throughout column A/rows
#    1-5) onto the handful of cells that currently carry one-off fonts, and
#    copy the "emphasis" style (style index 1, font size 16) that a couple
#    of new cells need -- do this BEFORE the new rows get their own values
#    so the format-only paste has an intact donor cell to copy from.
# ---------------------------------------------------------------------
$ws.Range("A2").Copy() | Out-Null
$normalTargets = @("A6","B7","D7","E7","B8","B9","C9","B10","C10")
foreach ($addr in $normalTargets) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# 2) Remove the stray empty, style-only cells left in column D (D8/D9/D10)
#    and D11 (which becomes entirely empty once cleared).
# ---------------------------------------------------------------------
$ws.Range("D8").Clear() | Out-Null
$ws.Range("D9").Clear() | Out-Null
$ws.Range("D10").Clear() | Out-Null
$ws.Range("D11").Clear() | Out-Null

# ---------------------------------------------------------------------
# 3) Rows 6-10 no longer need an explicit row height override -- let the
#    sheet's default height apply again.
# ---------------------------------------------------------------------
$ws.Rows(6).AutoFit() | Out-Null
$ws.Rows(7).AutoFit() | Out-Null
$ws.Rows(8).AutoFit() | Out-Null
$ws.Rows(9).AutoFit() | Out-Null
$ws.Rows(10).AutoFit() | Out-Null

# ---------------------------------------------------------------------
# 4) New Q&A rows (11-18) -- Reference Number / Message / Bot Response(/
#    User / Bot) continuing the existing table.
# ---------------------------------------------------------------------
$ws.Range("A11").Value2 = 1955
$ws.Range("B11").Value2 = "I saw messege in chat, I have completed testnet of 5 days but no role no nft received."
$ws.Range("C11").Value2 = "Have you change your discord name since you've completed the testnet tasks? What is the full ETH address?"
$ws.Range("A12").Value2 = 1954
$ws.Range("B12").Value2 = "the query is I have completed that five days task. but didn't get nft or role"
$ws.Range("C12").Value2 = "Please provide screenshots showing that you completed the steps. Also, make sure you have minted the NFT. NFT drops whitelist your address, but do not automatically deposit it."
$ws.Range("A13").Value2 = 1948
$ws.Range("B13").Value2 = "i can't claim my pioneer nft but i have claimed. my address 0x7fB87Ff912b81A9211fb4cca2445643702bf5D33"
$ws.Range("C13").Value2 = "Could you please clarify ... You can't claim but have claimed?"
$ws.Range("D13").Value2 = "i can't claim role hahah. and i can't mint them. so that's why am not eligble"
$ws.Range("E13").Value2 = "Without a tx id it's hard to know why these are failing, could you please provide that?"
$ws.Range("A14").Value2 = 1910
$ws.Range("B14").Value2 = "Dont get points. 0xBb41dd49254E8B9d631B835062392a460081734D."
$ws.Range("C14").Value2 = "Hello. If you have deposited recently, please refresh and check again later. More info regarding points can be found here. https://docs.primeprotocol.xyz/navigating-prime/prime-early-adopter-program"
$ws.Range("A15").Value2 = 1933
$ws.Range("B15").Value2 = "Hello. I noticed that when you repay the loan and withdraw the deposit, if you use the button 100%, then the entire amount is not withdrawn, the account balance is left, which will have to repay again and of course pay for gas."
$ws.Range("C15").Value2 = "Thank you for the information; I'll tag the team on this."
$ws.Range("A16").Value2 = 1942
$ws.Range("B16").Value2 = "0x9E3CeB09375f9d4922D6dc7Da3D006958298CDA8, id like to re-add my initial deposit - which was much more. really like how the P P team handled all this. im sure though you can understand my concerns about liquidity and rewards - i think you all will do great. lmk about the points"
$ws.Range("A17").Value2 = 1986
$ws.Range("B17").Value2 = "i deposit some glmr and dot, but I can't withdraw it. I found some errors in the execution of the contract, please help me to see if I can get my token back. my tx: 0xb6d674171cf570185a5ad3ddef47f853b71d4e923d139d0739d03c0bdc2fa319 0x0d54faeba552e9f5f0f4ad322800ed6f6539e4c4e74348171b09bfecb899831b"
$ws.Range("C17").Value2 = "Thank you for providing the transaction. I have tagged the team and they will continue from here."
$ws.Range("A18").Value2 = 1925
$ws.Range("B18").Value2 = "If I didn't get the daily connect posted yesterday before the zealy renewal time, doesn't this sprint count? After all, I'm shorted by 50xp for one daily connect, right?"
$ws.Range("C18").Value2 = "If you can't complete this task. You won't get XP"

# Rows 17 & 18 end up a bit taller than the sheet default (20pt instead of
# 16pt) and their "Message" column keeps the bigger/emphasis font used
# elsewhere in the sheet.
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$ws.Range("B7").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null

$ws.Rows(17).RowHeight = 20
$ws.Rows(18).RowHeight = 20

# ---------------------------------------------------------------------
# 5) Page setup + selection tweaks picked up by the saved view state.
# ---------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

$ws.Range("F16").Select() | Out-Null
